$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: bold text values
$ws.Range("A1").Value = "testando"
$ws.Range("B1").Value = "07/09/2025 17:16:35"

# Second row stays empty but still carries the bold style
$ws.Range("A1:B2").Font.Bold = $true

# Touch page setup so the worksheet emits a (empty) headerFooter element
$ws.PageSetup.CenterHeader = ""
